$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1830985915492958
$ws.Range("C2").Value = 0.5915492957746479
$ws.Range("J2").Value = 0.009389671361502348
$ws.Range("P2").Value = 0.1314553990610329
$ws.Range("S2").Value = 0.08450704225352113
$ws.Range("B3").Value = 0.01550387596899225
$ws.Range("C3").Value = 0.02325581395348837
$ws.Range("J3").Value = 0.02325581395348837
$ws.Range("P3").Value = 0.7131782945736435
$ws.Range("S3").Value = 0.2248062015503876
$ws.Range("J4").Value = 0.1290322580645161
$ws.Range("P4").Value = 0.5483870967741935
$ws.Range("S4").Value = 0.3225806451612903
$ws.Range("B6").Value = 0.05325443786982249
$ws.Range("D6").Value = 0.005917159763313609
$ws.Range("F6").Value = 0.03550295857988166
$ws.Range("J6").Value = 0.2603550295857988
$ws.Range("O6").Value = 0.02366863905325444
$ws.Range("Q6").Value = 0.1420118343195266
$ws.Range("R6").Value = 0.08284023668639054
$ws.Range("S6").Value = 0.3964497041420119
$ws.Range("B7").Value = 0.05555555555555555
$ws.Range("D7").Value = 0.02083333333333333
$ws.Range("F7").Value = 0.04166666666666666
$ws.Range("J7").Value = 0.2430555555555556
$ws.Range("Q7").Value = 0.1388888888888889
$ws.Range("R7").Value = 0.09722222222222222
$ws.Range("S7").Value = 0.4027777777777778
$ws.Range("B8").Value = 0.0584958217270195
$ws.Range("D8").Value = 0.01392757660167131
$ws.Range("F8").Value = 0.07520891364902507
$ws.Range("J8").Value = 0.116991643454039
$ws.Range("O8").Value = 0.02785515320334262
$ws.Range("Q8").Value = 0.1197771587743733
$ws.Range("R8").Value = 0.1197771587743733
$ws.Range("S8").Value = 0.467966573816156
$ws.Range("B9").Value = 0.07453416149068323
$ws.Range("D9").Value = 0.006211180124223602
$ws.Range("F9").Value = 0.05590062111801242
$ws.Range("J9").Value = 0.124223602484472
$ws.Range("O9").Value = 0.02484472049689441
$ws.Range("Q9").Value = 0.08695652173913043
$ws.Range("R9").Value = 0.1180124223602484
$ws.Range("S9").Value = 0.5093167701863354
$ws.Range("B10").Value = 0.102076124567474
$ws.Range("D10").Value = 0.01816608996539792
$ws.Range("E10").Value = 0.0008650519031141869
$ws.Range("F10").Value = 0.06055363321799308
$ws.Range("J10").Value = 0.1410034602076125
$ws.Range("O10").Value = 0.009515570934256055
$ws.Range("Q10").Value = 0.208477508650519
$ws.Range("R10").Value = 0.08477508650519031
$ws.Range("S10").Value = 0.3745674740484429
$ws.Range("G11").Value = 0.1428571428571428
$ws.Range("J11").Value = 0.1218487394957983
$ws.Range("K11").Value = 0.2310924369747899
$ws.Range("L11").Value = 0.4873949579831933
$ws.Range("S11").Value = 0.01680672268907563
$ws.Range("G12").Value = 0.7394957983193278
$ws.Range("J12").Value = 0.1848739495798319
$ws.Range("K12").Value = 0.01680672268907563
$ws.Range("L12").Value = 0.008403361344537815
$ws.Range("S12").Value = 0.05042016806722689
$ws.Range("G13").Value = 0.6578947368421053
$ws.Range("J13").Value = 0.3157894736842105
$ws.Range("S13").Value = 0.02631578947368421
$ws.Range("F15").Value = 0.02010050251256281
$ws.Range("H15").Value = 0.221105527638191
$ws.Range("I15").Value = 0.06030150753768844
$ws.Range("J15").Value = 0.3467336683417085
$ws.Range("K15").Value = 0.05527638190954774
$ws.Range("M15").Value = 0.01005025125628141
$ws.Range("O15").Value = 0.08040201005025126
$ws.Range("S15").Value = 0.2060301507537688
$ws.Range("F16").Value = 0.007575757575757576
$ws.Range("H16").Value = 0.1363636363636364
$ws.Range("I16").Value = 0.1136363636363636
$ws.Range("J16").Value = 0.4318181818181818
$ws.Range("K16").Value = 0.06818181818181818
$ws.Range("M16").Value = 0.01515151515151515
$ws.Range("O16").Value = 0.04545454545454546
$ws.Range("S16").Value = 0.1818181818181818
$ws.Range("F17").Value = 0.02064896755162242
$ws.Range("H17").Value = 0.1799410029498525
$ws.Range("I17").Value = 0.07079646017699115
$ws.Range("J17").Value = 0.4837758112094395
$ws.Range("K17").Value = 0.09144542772861357
$ws.Range("M17").Value = 0.02654867256637168
$ws.Range("O17").Value = 0.05014749262536873
$ws.Range("S17").Value = 0.07669616519174041
$ws.Range("F18").Value = 0.005376344086021506
$ws.Range("H18").Value = 0.1505376344086022
$ws.Range("I18").Value = 0.1075268817204301
$ws.Range("J18").Value = 0.4677419354838709
$ws.Range("K18").Value = 0.06451612903225806
$ws.Range("M18").Value = 0.01612903225806452
$ws.Range("O18").Value = 0.06451612903225806
$ws.Range("S18").Value = 0.1236559139784946
$ws.Range("F19").Value = 0.02102376599634369
$ws.Range("H19").Value = 0.1928702010968922
$ws.Range("I19").Value = 0.08409506398537477
$ws.Range("J19").Value = 0.3848263254113345
$ws.Range("K19").Value = 0.1005484460694698
$ws.Range("M19").Value = 0.02010968921389397
$ws.Range("N19").Value = 0.0009140767824497258
$ws.Range("O19").Value = 0.0850091407678245
$ws.Range("S19").Value = 0.1106032906764168
